$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Palmar" customer to full business name
$ws.Range("A27").Value = "PALMAR MEXICAN RESTAURANT"

# Fix C27 prospect/salesperson code (013 -> 023)
$ws.Range("C27").Value = "023"

# Update several Last Invoice Date values (file-based timestamps)
$ws.Range("D3").Value = 45868
$ws.Range("D10").Value = 45874
$ws.Range("D11").Value = 45849
$ws.Range("D12").Value = 45845
$ws.Range("D30").Value = 45848

# Remove the now-unused "Last Invoice Date" helper column E entirely
$ws.Range("E2:E30").Clear()
